$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.342.83"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "2.446.65"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.23"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.35"
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.532"
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").Value = "2.444.60"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  -3.29%  "
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.345"
$ws.Range("E13").Value = "  -2.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.58"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("E15").Value = "  -2.96%  "
$ws.Range("D16").Value = "2.887.48"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "62.100.53"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "2.445.83"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.93"
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.16"
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "330.73"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("E23").Value = "  -3.71%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.96"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.51"
$ws.Range("E26").Value = "  +7.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "624.78"
$ws.Range("E27").Value = "  +1.90%  "
$ws.Range("D28").Value = "2.566.78"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("E29").Value = "  -5.78%  "
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("E31").Value = "  -3.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.04"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("E33").Value = "  +0.74%  "
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("E35").Value = "  -4.59%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("E37").Value = "  -5.26%  "
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "149.72"
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.30"
$ws.Range("E40").Value = "  -1.66%  "
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("E42").Value = "  -1.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.55"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -4.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.89"
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.602"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("E50").Value = "  -6.94%  "
$ws.Range("D51").Value = "0.0₆0238"
$ws.Range("E51").Value = "  +8.99%  "
